# Insert a new data row at row 23 (shifting the existing rows 23:103 down
# to 24:104) and populate it with the new "Santina" price entry, matching
# the weekly refresh described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 23 and below down by one row.
$ws.Rows.Item(23).Insert()

# Populate the freshly inserted row 23 with the new record.
$ws.Cells.Item(23, 1).Value = 7
$ws.Cells.Item(23, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(23, 3).Value = "Ñuble"
$ws.Cells.Item(23, 4).Value = 44592
$ws.Cells.Item(23, 5).Value = 16
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100103
$ws.Cells.Item(23, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(23, 9).Value = 100103001
$ws.Cells.Item(23, 10).Value = "Cereza"
$ws.Cells.Item(23, 11).Value = "Santina"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 120
$ws.Cells.Item(23, 14).Value = 5500
$ws.Cells.Item(23, 15).Value = 6000
$ws.Cells.Item(23, 16).Value = 5750
$ws.Cells.Item(23, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(23, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(23, 19).Value = 575
$ws.Cells.Item(23, 20).Value = 10
